# Auto-generated edit script: update "想去人数" (column F) counts
# per the diff, across four worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 285
$ws.Range("F6").Value = 1747
$ws.Range("F7").Value = 2097
$ws.Range("F8").Value = 1471
$ws.Range("F9").Value = 1304
$ws.Range("F10").Value = 2986
$ws.Range("F11").Value = 1879
$ws.Range("F12").Value = 1317
$ws.Range("F13").Value = 1732
$ws.Range("F15").Value = 160
$ws.Range("F16").Value = 600
$ws.Range("F18").Value = 2025
$ws.Range("F22").Value = 1749
$ws.Range("F23").Value = 659
$ws.Range("F24").Value = 4944
$ws.Range("F34").Value = 2922
$ws.Range("F37").Value = 3444
$ws.Range("F43").Value = 253

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 148579
$ws.Range("F8").Value = 148579
$ws.Range("F19").Value = 288
$ws.Range("F30").Value = 141
$ws.Range("F36").Value = 275
$ws.Range("F39").Value = 210

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F10").Value = 314
$ws.Range("F11").Value = 2500
$ws.Range("F12").Value = 79
$ws.Range("F13").Value = 914

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 314
$ws.Range("F8").Value = 2500
$ws.Range("F9").Value = 285
$ws.Range("F12").Value = 148579
$ws.Range("F13").Value = 1747
$ws.Range("F14").Value = 2097
$ws.Range("F15").Value = 1471
$ws.Range("F16").Value = 1304
$ws.Range("F17").Value = 2986
$ws.Range("F18").Value = 1879
$ws.Range("F19").Value = 1317
$ws.Range("F21").Value = 1732
$ws.Range("F22").Value = 160
$ws.Range("F23").Value = 600
$ws.Range("F24").Value = 288
$ws.Range("F26").Value = 2025
$ws.Range("F30").Value = 1749
$ws.Range("F31").Value = 4944
$ws.Range("F37").Value = 914
$ws.Range("F43").Value = 2922
$ws.Range("F46").Value = 3444
$ws.Range("F47").Value = 275
$ws.Range("F51").Value = 254
